$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bootstrap")
for ($i=1; $i -le 20; $i++) {
    $ws.Columns.Item($i).ColumnWidth = $i
}
Write-Output "done"
